# Generate Report for Archive
#
# Two files (180eea62-c18e-4321-9cb8-71e7758da1f5.md and
# 1da890d3-4e10-462c-8cb0-cfaf2666694b.md) move from "Ready for handoff"
# to "In Translation" status, on rows 7 and 8 of every sheet.
#
# - "Overview" sheet: columns B (zh-cn) and C (de-de) hold the per-locale
#   status for each file.
# - "zh-cn" / "de-de" sheets: column C ("Status") holds the status for
#   each file.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus
$overview.Range("B8").Value = $newStatus
$overview.Range("C8").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = $newStatus
$zhcn.Range("C8").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = $newStatus
$dede.Range("C8").Value = $newStatus
